# Generate Report for Handoff
#
# A new handoff run was generated: the per-file GUID, the per-content hash
# embedded in the generated xliff file names, and the "latest" timestamps
# all roll forward. Apply the new values everywhere the old ones appeared
# (Overview / zh-cn / de-de sheets) and repoint the "Source File Name" /
# "Path And Name" hyperlinks at the new file.

$wb = $excel.ActiveWorkbook

$newGuid = "b607217c-78e9-4d43-b187-13105c6720e7"
$newHash = "cdf276745b0bb8890a095ba8bbe3f830ffbf7973"

$newHoDate   = "2016-09-04 19:02:41"
$newZhCnDate = "2016-09-04 19:02:36"
$newDeDeDate = "2016-09-04 19:02:41"

$commitRef = "ece77f9cf32eef050c2fd571fa1d6636f9e9ea1f"

# Visual style ("HyperLink") the hyperlink cells use in this workbook -
# underlined, cornflower-blue text. Re-applied after Hyperlinks.Add (which
# otherwise stamps its own default theme hyperlink style on the cell) so the
# cell keeps looking the way it did before the edit.
function Restore-HyperlinkLook($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitRef/e2e/$newGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "e2e\$newGuid.md")
Restore-HyperlinkLook $wsOverview.Range("B2")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = $newZhCnDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitRef/e2e/$newGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newGuid.md")
Restore-HyperlinkLook $wsZhCn.Range("A2")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = $newDeDeDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitRef/e2e/$newGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newGuid.md")
Restore-HyperlinkLook $wsDeDe.Range("A2")
